$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2231.7368
$ws.Range("I70").Value = 1225.5
$ws.Range("J70").Value = 2500.0667
$ws.Range("K70").Value = 3676.5
$ws.Range("L70").Value = 7500.2001
$ws.Range("M70").Value = -3406.5
$ws.Range("N70").Value = -8040.2001

$ws.Range("H73").Value = 2231.7368
$ws.Range("I73").Value = 1225.5
$ws.Range("J73").Value = 2500.0667
$ws.Range("K73").Value = 3676.5
$ws.Range("L73").Value = 7500.2001
$ws.Range("M73").Value = -2740.5
$ws.Range("N73").Value = -9372.2001

$ws.Range("H97").Value = 2414.9167
$ws.Range("J97").Value = 2414.9167
$ws.Range("L97").Value = 7244.750100000001
$ws.Range("N97").Value = -8236.750100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1093412.8
$ws.Range("I32").Value = 1210017.5
$ws.Range("J32").Value = 18058.777
$ws.Range("K32").Value = 1210017.5
$ws.Range("L32").Value = 18058.777
$ws.Range("M32").Value = -1209730.5
$ws.Range("N32").Value = -18632.777

$ws.Range("H88").Value = 23600.555
$ws.Range("J88").Value = 34733.168
$ws.Range("L88").Value = 34733.168
$ws.Range("N88").Value = -35545.168

$ws.Range("H91").Value = 23600.555
$ws.Range("J91").Value = 34733.168
$ws.Range("L91").Value = 34733.168
$ws.Range("N91").Value = -37541.168

$ws.Range("H119").Value = 31000
$ws.Range("J119").Value = 31000
$ws.Range("L119").Value = 31000
$ws.Range("N119").Value = -40676

$ws.Range("H132").Value = 1969617.8
$ws.Range("I132").Value = 2269993.2
$ws.Range("J132").Value = 843209.4399999999
$ws.Range("K132").Value = 6809979.600000001
$ws.Range("L132").Value = 2529628.32
$ws.Range("M132").Value = -6807449.600000001
$ws.Range("N132").Value = -2534688.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 822.7143
$ws.Range("I7").Value = 1015.2727
$ws.Range("J7").Value = 116.666664
$ws.Range("K7").Value = 1015.2727
$ws.Range("L7").Value = 116.666664
$ws.Range("M7").Value = -902.2727
$ws.Range("N7").Value = -342.666664

$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H58").Value = 1156.6786
$ws.Range("I58").Value = 779.17645
$ws.Range("J58").Value = 1740.091
$ws.Range("K58").Value = 779.17645
$ws.Range("L58").Value = 1740.091
$ws.Range("M58").Value = -576.17645
$ws.Range("N58").Value = -2146.091

$ws.Range("H136").Value = 1156.6786
$ws.Range("I136").Value = 779.17645
$ws.Range("J136").Value = 1740.091
$ws.Range("K136").Value = 2337.52935
$ws.Range("L136").Value = 5220.272999999999
$ws.Range("M136").Value = 212.4706499999998
$ws.Range("N136").Value = -10320.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7233.8
$ws.Range("I5").Value = 516.8889
$ws.Range("J5").Value = 17309.166
$ws.Range("K5").Value = 1550.6667
$ws.Range("L5").Value = 51927.49800000001
$ws.Range("M5").Value = -1438.6667
$ws.Range("N5").Value = -52151.49800000001

$ws.Range("H122").Value = 632.5185
$ws.Range("I122").Value = 412
$ws.Range("J122").Value = 709.7
$ws.Range("K122").Value = 3708
$ws.Range("L122").Value = 6387.3
$ws.Range("M122").Value = -1258
$ws.Range("N122").Value = -11287.3

$ws.Range("H135").Value = 7233.8
$ws.Range("I135").Value = 516.8889
$ws.Range("J135").Value = 17309.166
$ws.Range("K135").Value = 4652.0001
$ws.Range("L135").Value = 155782.494
$ws.Range("M135").Value = -2117.0001
$ws.Range("N135").Value = -160852.494

$ws.Range("H137").Value = 8969518
$ws.Range("I137").Value = 15152419
$ws.Range("J137").Value = 5568923
$ws.Range("K137").Value = 45457257
$ws.Range("L137").Value = 16706769
$ws.Range("M137").Value = -45452157
$ws.Range("N137").Value = -16716969

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 876.9375
$ws.Range("I122").Value = 809.46155
$ws.Range("J122").Value = 1169.3334
$ws.Range("K122").Value = 2428.38465
$ws.Range("L122").Value = 3508.0002
$ws.Range("M122").Value = 21.61535000000003
$ws.Range("N122").Value = -8408.0002

$ws.Range("H132").Value = 34500.902
$ws.Range("I132").Value = 1911.5333
$ws.Range("J132").Value = 65053.438
$ws.Range("K132").Value = 5734.5999
$ws.Range("L132").Value = 195160.314
$ws.Range("M132").Value = -3204.5999
$ws.Range("N132").Value = -200220.314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1993.2273
$ws.Range("I61").Value = 1786.909
$ws.Range("J61").Value = 2199.5454
$ws.Range("K61").Value = 1786.909
$ws.Range("L61").Value = 2199.5454
$ws.Range("M61").Value = -1584.909
$ws.Range("N61").Value = -2603.5454

$ws.Range("H113").Value = 1993.2273
$ws.Range("I113").Value = 1786.909
$ws.Range("J113").Value = 2199.5454
$ws.Range("K113").Value = 1786.909
$ws.Range("L113").Value = 2199.5454
$ws.Range("M113").Value = 383.0909999999999
$ws.Range("N113").Value = -6539.5454

$ws.Range("H119").Value = 21000
$ws.Range("J119").Value = 21000
$ws.Range("L119").Value = 21000
$ws.Range("N119").Value = -30676

$ws.Range("H132").Value = 182132.83
$ws.Range("I132").Value = 48449.465
$ws.Range("K132").Value = 145348.395
$ws.Range("M132").Value = -142818.395

$ws.Range("H136").Value = 239697.78
$ws.Range("I136").Value = 313796.78
$ws.Range("J136").Value = 2581
$ws.Range("K136").Value = 941390.3400000001
$ws.Range("L136").Value = 7743
$ws.Range("M136").Value = -938840.3400000001
$ws.Range("N136").Value = -12843

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3442.8096
$ws.Range("I96").Value = 1801
$ws.Range("J96").Value = 3716.4443
$ws.Range("K96").Value = 1801
$ws.Range("L96").Value = 3716.4443
$ws.Range("M96").Value = -428
$ws.Range("N96").Value = -6462.4443

$ws.Range("H119").Value = 3000
$ws.Range("J119").Value = 3000
$ws.Range("L119").Value = 3000
$ws.Range("N119").Value = -12676

$ws.Range("H132").Value = 1969.9359
$ws.Range("I132").Value = 472.36667
$ws.Range("J132").Value = 6961.8335
$ws.Range("K132").Value = 1417.10001
$ws.Range("L132").Value = 20885.5005
$ws.Range("M132").Value = 1112.89999
$ws.Range("N132").Value = -25945.5005

$ws.Range("H136").Value = 1164055.4
$ws.Range("I136").Value = 1191180.4
$ws.Range("J136").Value = 1001305
$ws.Range("K136").Value = 3573541.2
$ws.Range("L136").Value = 3003915
$ws.Range("M136").Value = -3570991.2
$ws.Range("N136").Value = -3009015
